$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 data values
$ws.Range("B2").Value = 16.656345414395474
$ws.Range("C2").Value = 13.102169094515716
$ws.Range("D2").Value = 13.230320720899556
$ws.Range("E2").Value = 0.50889455340507084

# Update row 3 data values
$ws.Range("B3").Value = 33.916586186595538
$ws.Range("C3").Value = 4.0443387243297195
$ws.Range("D3").Value = 2.768088959534424
$ws.Range("E3").Value = 3.0261117330107936

# Update the selection on the sheet to match the new active range
$ws.Range("B1:E3").Select()
